# Trade #76 closed at 2026-02-17 08:58:14 - unknown UNKNOWN +0.000%
# Updates Summary + Strategy Status aggregate stats and appends the new
# trade row to both the "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# Helper: write a value as literal TEXT (not auto-parsed into a date/number)
# by forcing a Text number format on entry, then resetting the cell style
# back to "Normal" afterwards so no stray formatting is left behind.
function Set-TextValue($ws, $cellRef, $val) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range($cellRef).Style = "Normal"
}

# ---- Summary sheet -------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 76
$summary.Range("B9").Value = 40.79

# ---- Strategy Status sheet ------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 76
$status.Range("G4").Value = 40.79

# ---- Append new trade row (#76) to "All Trades" and "MarketMaking" --
$rowNum = 77

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A$rowNum").Value = 76
    Set-TextValue $ws "B$rowNum" "2026-02-17"
    Set-TextValue $ws "C$rowNum" "08:58:08"
    Set-TextValue $ws "D$rowNum" "MarketMaking"
    Set-TextValue $ws "E$rowNum" "UP"
    $ws.Range("F$rowNum").Value = 0.15
    $ws.Range("G$rowNum").Value = 0.15
    Set-TextValue $ws "H$rowNum" "CLOSED"
    $ws.Range("I$rowNum").Value = 0
    $ws.Range("J$rowNum").Value = 0
    $ws.Range("K$rowNum").Value = 100.52
    $ws.Range("L$rowNum").Value = 0
    $ws.Range("M$rowNum").Value = 0
    $ws.Range("N$rowNum").Value = 0.6
    Set-TextValue $ws "O$rowNum" "Normal spread capture: 19600 bps"
    Set-TextValue $ws "P$rowNum" "early_exit"
    $ws.Range("Q$rowNum").Value = 0.13
}
